# Update the two listing rows with the new item info, and clear out the
# placeholder rows (4-10) that previously held quote-prefixed "blank" text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: first item
$ws.Range("A2").Value = "DSE 校園經濟 5蚊一本 有5本"
$ws.Range("B2").Value = "Brand new"
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = "新淨"

# Row 3: second item
$ws.Range("A3").Value = "DSE 公社刊物 5 蚊一本 有3本"
$ws.Range("B3").Value = "Brand new"
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = "新淨"

# Row 3 now wraps like row 2, so it grows to the same height.
$ws.Rows.Item(3).RowHeight = 19.5

# Rows 4-10 were placeholder rows carrying quote-prefixed empty text; clear
# them back to genuinely empty cells.
for ($r = 4; $r -le 10; $r++) {
    $ws.Cells.Item($r, 1).ClearContents()
    $ws.Cells.Item($r, 2).ClearContents()
    $ws.Cells.Item($r, 4).ClearContents()
}

# Re-apply the (non quote-prefixed) formatting used by row 2 to rows 4-10,
# replacing the old quote-prefix styles.
$ws.Range("A2:B2").Copy()
$ws.Range("A4:B10").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D4:D10").PasteSpecial(-4122)
